$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 21826.555
$ws.Range("I21").Value = 22919.572
$ws.Range("J21").Value = 18001
$ws.Range("K21").Value = 22919.572
$ws.Range("L21").Value = 18001
$ws.Range("M21").Value = -22451.572
$ws.Range("N21").Value = -18937

$ws.Range("H23").Value = 21826.555
$ws.Range("I23").Value = 22919.572
$ws.Range("J23").Value = 18001
$ws.Range("K23").Value = 22919.572
$ws.Range("L23").Value = 18001
$ws.Range("M23").Value = -22685.572
$ws.Range("N23").Value = -18469

$ws.Range("H40").Value = 2180.9524
$ws.Range("I40").Value = 1856.6666
$ws.Range("J40").Value = 2310.6667
$ws.Range("K40").Value = 1856.6666
$ws.Range("L40").Value = 2310.6667
$ws.Range("M40").Value = -1681.6666
$ws.Range("N40").Value = -2660.6667

$ws.Range("H58").Value = 657
$ws.Range("I58").Value = 142.5
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 427.5
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -277.5
$ws.Range("N58").Value = -3300

$ws.Range("H64").Value = 100002696
$ws.Range("I64").Value = 166668750
$ws.Range("J64").Value = 3615
$ws.Range("K64").Value = 166668750
$ws.Range("L64").Value = 3615
$ws.Range("M64").Value = -166668502
$ws.Range("N64").Value = -4111

$ws.Range("H67").Value = 100002696
$ws.Range("I67").Value = 166668750
$ws.Range("J67").Value = 3615
$ws.Range("K67").Value = 166668750
$ws.Range("L67").Value = 3615
$ws.Range("M67").Value = -166667892
$ws.Range("N67").Value = -5331

$ws.Range("H106").Value = 1396
$ws.Range("I106").Value = 1285.6
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 1285.6
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -654.5999999999999
$ws.Range("N106").Value = -3762

$ws.Range("H107").Value = 1005
$ws.Range("I107").Value = 1005
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1005
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 915

$ws.Range("H112").Value = 1227.5
$ws.Range("I112").Value = 933.3333
$ws.Range("J112").Value = 2110
$ws.Range("K112").Value = 2799.9999
$ws.Range("L112").Value = 6330
$ws.Range("M112").Value = -1691.9999
$ws.Range("N112").Value = -8546

$ws.Range("H118").Value = 1267.8125
$ws.Range("I118").Value = 481.42856
$ws.Range("J118").Value = 1488
$ws.Range("K118").Value = 1444.28568
$ws.Range("L118").Value = 4464
$ws.Range("M118").Value = 212.71432
$ws.Range("N118").Value = -7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4614.314
$ws.Range("I32").Value = 5230.075
$ws.Range("J32").Value = 2375.182
$ws.Range("K32").Value = 5230.075
$ws.Range("L32").Value = 2375.182
$ws.Range("M32").Value = -4943.075
$ws.Range("N32").Value = -2949.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1312.8948
$ws.Range("I94").Value = 775.4167
$ws.Range("J94").Value = 2234.2856
$ws.Range("K94").Value = 775.4167
$ws.Range("L94").Value = 2234.2856
$ws.Range("M94").Value = -324.4167
$ws.Range("N94").Value = -3136.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 662.3542
$ws.Range("I5").Value = 394.2143
$ws.Range("J5").Value = 1037.75
$ws.Range("K5").Value = 1182.6429
$ws.Range("L5").Value = 3113.25
$ws.Range("M5").Value = -1070.6429
$ws.Range("N5").Value = -3337.25

$ws.Range("H20").Value = 1280
$ws.Range("I20").Value = 975
$ws.Range("J20").Value = 2500
$ws.Range("K20").Value = 2925
$ws.Range("L20").Value = 7500
$ws.Range("M20").Value = -2698
$ws.Range("N20").Value = -7954

$ws.Range("H21").Value = 798.0909
$ws.Range("I21").Value = 474.625
$ws.Range("J21").Value = 1660.6666
$ws.Range("K21").Value = 1423.875
$ws.Range("L21").Value = 4981.9998
$ws.Range("M21").Value = -1250.875
$ws.Range("N21").Value = -5327.9998

$ws.Range("H22").Value = 20834708
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 22223622
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 66670866
$ws.Range("M22").Value = -2831
$ws.Range("N22").Value = -66671204

$ws.Range("H27").Value = 20834708
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 22223622
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 66670866
$ws.Range("M27").Value = -2898
$ws.Range("N27").Value = -66671070

$ws.Range("H34").Value = 497.77777
$ws.Range("I34").Value = 354.2857
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1062.8571
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -978.8571000000002
$ws.Range("N34").Value = -3168

$ws.Range("H40").Value = 247.5
$ws.Range("I40").Value = 240
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 960
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -891
$ws.Range("N40").Value = -1338

$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1409

$ws.Range("H57").Value = 4387
$ws.Range("I57").Value = 978.3333
$ws.Range("J57").Value = 9500
$ws.Range("K57").Value = 2934.9999
$ws.Range("L57").Value = 28500
$ws.Range("M57").Value = -2375.9999
$ws.Range("N57").Value = -29618

$ws.Range("H58").Value = 166667170
$ws.Range("I58").Value = 985
$ws.Range("J58").Value = 333333340
$ws.Range("K58").Value = 2955
$ws.Range("L58").Value = 1000000020
$ws.Range("M58").Value = -2827
$ws.Range("N58").Value = -1000000276

$ws.Range("H102").Value = 2780.6155
$ws.Range("I102").Value = 890
$ws.Range("J102").Value = 3347.8
$ws.Range("K102").Value = 2670
$ws.Range("L102").Value = 10043.4
$ws.Range("M102").Value = -236
$ws.Range("N102").Value = -14911.4

$ws.Range("H122").Value = 1206.125
$ws.Range("I122").Value = 630.4
$ws.Range("J122").Value = 2165.6667
$ws.Range("K122").Value = 5673.599999999999
$ws.Range("L122").Value = 19491.0003
$ws.Range("M122").Value = -3223.599999999999
$ws.Range("N122").Value = -24391.0003

$ws.Range("H131").Value = 651.7
$ws.Range("I131").Value = 284.1591
$ws.Range("J131").Value = 940.4820999999999
$ws.Range("K131").Value = 852.4773
$ws.Range("L131").Value = 2821.4463
$ws.Range("M131").Value = 4187.5227
$ws.Range("N131").Value = -12901.4463

$ws.Range("H135").Value = 662.3542
$ws.Range("I135").Value = 394.2143
$ws.Range("J135").Value = 1037.75
$ws.Range("K135").Value = 3547.9287
$ws.Range("L135").Value = 9339.75
$ws.Range("M135").Value = -1012.9287
$ws.Range("N135").Value = -14409.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 48000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 48000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 48000
$ws.Range("N88").Value = -48902

$ws.Range("H91").Value = 48000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 48000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 48000
$ws.Range("N91").Value = -51120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2998
$ws.Range("I93").Value = 2998
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2998
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -1750
